# Apply edits described by the commit:
# "Created Function for Gaussian Quadrature Scheme, and exported it to the
#  Averaged Intensities files."
#
# Changes:
#  1. Rename the worksheet from "GammaFiber2F-HW15.xpc" to "GammaFiber2F"
#  2. Minor floating point recalculation updates on row 15 (C15, F15, H15)
#  3. Append a new data row (row 16) with averaged intensity results

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename worksheet
$ws.Name = "GammaFiber2F"

# 2. Update slightly-recalculated values on row 15
$ws.Range("C15").Value = 0.9783807829387059
$ws.Range("F15").Value = 0.9783807829387059
$ws.Range("H15").Value = 0.7886986867608025

# 3. Append new row 16 with averaged intensity data
# Copy the formatting from A15 (bold, centered, bordered) onto A16, matching
# the style used for the other "HKL index" cells in column A.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.36567837976243
$ws.Range("D16").Value = 0.3828158769166418
$ws.Range("E16").Value = 1.050977968836757
$ws.Range("F16").Value = 1.36567837976243
$ws.Range("G16").Value = 0.6858759666533483
$ws.Range("H16").Value = 1.133475496403593
$ws.Range("I16").Value = 1.13367014862962
$ws.Range("J16").Value = 0.3828158769166418
$ws.Range("K16").Value = 0.7168969228766995
$ws.Range("L16").Value = 1.041287651319565
$ws.Range("M16").Value = 0.958748972867065
